$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '45.334.99'
$ws.Range("E2").Value = '  +2.67%  '
$ws.Range("D3").Value = '2.428.30'
$ws.Range("E3").Value = '  +0.03%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '318.48'
$ws.Range("E5").Value = '  +3.40%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '102.84'
$ws.Range("E6").Value = '  +2.25%  '
$ws.Range("E7").Value = '  +0.54%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("E8").Value = '  -0.09%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.529'
$ws.Range("E9").Value = '  +5.78%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '35.64'
$ws.Range("E10").Value = '  +0.76%  '
$ws.Range("E11").Value = '  +0.02%  '
$ws.Range("E12").Value = '  -2.06%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '18.23'
$ws.Range("E13").Value = '  -2.57%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '7.10'
$ws.Range("E14").Value = '  +2.35%  '
$ws.Range("D15").Value = '2.807.50'
$ws.Range("E15").Value = '  -0.03%  '
$ws.Range("D16").Value = '2.437.83'
$ws.Range("E16").Value = '  +1.61%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.843'
$ws.Range("E17").Value = '  +1.03%  '
$ws.Range("D18").Value = '45.223.36'
$ws.Range("E18").Value = '  +2.39%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.24'
$ws.Range("E19").Value = '  -0.75%  '
$ws.Range("E20").Value = '  -1.69%  '
$ws.Range("E21").Value = '  +1.54%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '68.90'
$ws.Range("E22").Value = '  +0.42%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '244.48'
$ws.Range("E23").Value = '  +1.63%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.27'
$ws.Range("E24").Value = '  -1.42%  '
$ws.Range("E25").Value = '  +0.24%  '
$ws.Range("E26").Value = '  +0.05%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '25.74'
$ws.Range("E27").Value = '  +1.85%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.19'
$ws.Range("E28").Value = '  -6.78%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.61'
$ws.Range("E29").Value = '  -0.11%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '49.43'
$ws.Range("E30").Value = '  +2.91%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '33.05'
$ws.Range("E31").Value = '  -1.16%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.126'
$ws.Range("E32").Value = '  +6.16%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '20.03'
$ws.Range("E33").Value = '  +7.09%  '
$ws.Range("E34").Value = '  +0.57%  '
$ws.Range("E35").Value = '  +0.21%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0765'
$ws.Range("E36").Value = '  +0.06%  '
$ws.Range("E37").Value = '  -2.23%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.45'
$ws.Range("E38").Value = '  -1.75%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.89'
$ws.Range("E39").Value = '  -1.52%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '126.31'
$ws.Range("E40").Value = '  -3.01%  '
$ws.Range("E41").Value = '  +0.83%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.20'
$ws.Range("E42").Value = '  -3.76%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '20.72'
$ws.Range("E43").Value = '  -1.47%  '
$ws.Range("E44").Value = '  +0.74%  '
$ws.Range("D45").Value = '1.939.08'
$ws.Range("E45").Value = '  -1.11%  '
$ws.Range("E46").Value = '  -2.68%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.94'
$ws.Range("E47").Value = '  +1.98%  '
$ws.Range("E48").Value = '  +9.47%  '
$ws.Range("E49").Value = '  -3.17%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '76.96'
$ws.Range("E50").Value = '  +4.73%  '

# Row 51: THORChain -> MultiversX
$ws.Range("B51").Value = 'MultiversX'
$ws.Range("C51").Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '54.15'
$ws.Range("E51").Value = '  +1.22%  '
